$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M30").Value = 61.75
$ws1.Range("H45").Value = 71.09999999999999
$ws1.Range("I45").Value = 100.8
$ws1.Range("M46").Value = 180.61
$ws1.Range("H57").Value = "1 de 55"
$ws1.Range("I57").Value = "2 de 55"
$ws1.Range("M57").Value = "3 de 55"

$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F30").Value = 61.75
$ws2.Range("F45").Value = 171.9
$ws2.Range("F46").Value = 180.61
$ws2.Range("F57").Value = 11823.47

$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 10165.82
$ws3.Range("E3").Value = 4718.059999999999
$ws3.Range("F3").Value = 0.6830087315941811

$ws3.Range("D7").Value = 71.09999999999999
$ws3.Range("E7").Value = 3128.9
$ws3.Range("F7").Value = 0.02221875

$ws3.Range("D8").Value = -52.8
$ws3.Range("E8").Value = 1052.8
$ws3.Range("F8").Value = -0.0528

$ws3.Range("D15").Value = 3046.66
$ws3.Range("E15").Value = 17643.34
$ws3.Range("F15").Value = 0.147252779120348

$ws3.Range("D16").Value = 1047.06
$ws3.Range("E16").Value = 57674.17000000001
$ws3.Range("F16").Value = 0.01783102976555498

$ws3.Range("D19").Value = 15023.6
$ws3.Range("E19").Value = 94845.15000000001
$ws3.Range("F19").Value = 0.1367413390977871

$ws3.Columns.Item(4).ColumnWidth = 13.166666666666666
